# Applies the "Fixed attack sliding bug" update to the Sprint 1 planning sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task #6 "Player class" - real time spent (H13) corrected from 300 to 340 minutes.
$ws.Range("H13").Value2 = 340

# Task #12 "Main Menu" (row 19) now has real time recorded, completion percentage,
# and a comment explaining remaining work.
$ws.Range("H19").Value2 = 150
$ws.Range("K19").Value2 = "75%"
$ws.Range("L19").Value2 = "Autres options comme Settings et Load Game à venir, plus long car changement concept"

# Task #8 "UI Manager (switch between menus)" comment updated.
$ws.Range("L15").Value2 = "Pas d'autre level pour le moment "

# Header row grew taller to accommodate the updated content.
$ws.Rows(3).RowHeight = 63

# Leave the cursor on the cell that was actually edited (H13), matching the
# final selection state of the workbook.
$ws.Range("H13").Select()
